$d = $word.ActiveDocument

# New grid-column / cell widths (dxa -> points, 20 dxa = 1 pt)
$colWidths = @(69.05, 84.9, 219.7, 70.8, 92.1, 106.25, 64.2)

for ($ti = 1; $ti -le $d.Tables.Count; $ti++) {
    $t = $d.Tables.Item($ti)

    # w:tblInd  0 -> -5 dxa  (-0.25 pt)
    $t.Rows.LeftIndent = -0.25

    # w:tblCellMar left  108 -> 103 dxa (5.4 -> 5.15 pt)
    $t.LeftPadding = 5.15

    # w:gridCol / w:tcW widths
    for ($ci = 1; $ci -le $t.Columns.Count; $ci++) {
        $t.Columns.Item($ci).Width = $colWidths[$ci - 1]
    }

    # w:tcMar left 108 -> 103 dxa (5.4 -> 5.15 pt) on every cell
    foreach ($cell in $t.Range.Cells) {
        $cell.LeftPadding = 5.15
    }

    # Header text "POSIÇÃO SAL" -> "POSIÇÃO MJC"
    $cell = $t.Cell(1, 4)
    $cell.Range.Find.Execute("SAL", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "MJC", 2) | Out-Null
}
